# Append new trade row (#33) to the "leadlag" sheet, closing out the diff's
# row 29 addition (dimension grows from A1:N28 to A1:N29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 29

# Numeric columns
$ws.Cells.Item($row, 1).Value  = 33                  # A: Trade #
$ws.Cells.Item($row, 6).Value  = 68849.49000000001   # F: Entry Price
$ws.Cells.Item($row, 9).Value  = 0                   # I: P&L %
$ws.Cells.Item($row, 10).Value = 0                   # J: P&L $
$ws.Cells.Item($row, 11).Value = 0.75                # K: Confidence
$ws.Cells.Item($row, 14).Value = 0                   # N: Duration (min)

# Text columns - leading apostrophe forces literal text (avoids
# date/time/number auto-conversion), matching the source inlineStr cells.
$ws.Cells.Item($row, 2).Value  = "'2026-02-16"                        # B: Date
$ws.Cells.Item($row, 3).Value  = "'21:28:28"                          # C: Time
$ws.Cells.Item($row, 4).Value  = "leadlag"                            # D: Strategy
$ws.Cells.Item($row, 5).Value  = "DOWN"                               # E: Side
$ws.Cells.Item($row, 8).Value  = "OPEN"                               # H: Status
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.160% move"  # L: Entry Reason

# Exit Price (G) and Exit Reason (M) stay blank-but-text, same as every
# other still-OPEN trade row (e.g. row 28): an empty quoted string keeps
# the cell as empty text instead of clearing it to a truly blank cell.
$ws.Cells.Item($row, 7).Value  = "'"   # G: Exit Price
$ws.Cells.Item($row, 13).Value = "'"   # M: Exit Reason
